# Update "想去人数" (F column) figures for several events.
# The same events are listed both on the "展览" sheet and on the
# aggregated "全部类型" sheet, so each value needs to be updated in
# both places.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# row on 展览 -> row on 全部类型 -> new value
$updates = @(
    @{ ExpoRow = 2;  AllRow = 2;  Value = 456 },  # W·A第五人格同人only2.0
    @{ ExpoRow = 3;  AllRow = 3;  Value = 38 },   # 之心城购物中心-2024漫趣地带嘉年华
    @{ ExpoRow = 4;  AllRow = 7;  Value = 37 },   # 崩坏同人only 爱莉希雅同人生日会
    @{ ExpoRow = 6;  AllRow = 9;  Value = 5292 }, # 第九届环形宇宙动漫游戏嘉年华
    @{ ExpoRow = 7;  AllRow = 10; Value = 179 },  # MAX特摄同人only2.0
    @{ ExpoRow = 8;  AllRow = 11; Value = 101 },  # 九号幻想动漫游戏嘉年华
    @{ ExpoRow = 9;  AllRow = 12; Value = 104 },  # 风月引代号鸢同人only
    @{ ExpoRow = 10; AllRow = 14; Value = 362 },  # 心动恋章·冬日序国乙&代号鸢同人only
    @{ ExpoRow = 11; AllRow = 15; Value = 51 }    # 星光国潮动漫游戏嘉年华
)

foreach ($u in $updates) {
    $wsExpo.Cells.Item($u.ExpoRow, 6).Value = $u.Value
    $wsAll.Cells.Item($u.AllRow, 6).Value = $u.Value
}
